# Added Flow vs R1L to the cell data modeled by tissue slice code
#
# This mirrors the existing "Kpl" summary table (rows 22-24, built off
# column B) by adding a matching summary table for "Flow_Lac" (rows
# 38-40, built off column F): a header row of the four sample labels,
# a row of per-group averages, and a row of per-group standard errors.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38: repeat the group labels (same shared strings as row 22) above
# the new Flow_Lac summary block.
$ws.Range("G38").Value = "HK-2"
$ws.Range("H38").Value = "UMRC6"
$ws.Range("I38").Value = "UOK262"
$ws.Range("J38").Value = "UOK + DIDS"

# Row 39: row label + average of column F (Flow_Lac) per group.
$ws.Range("F39").Value = "Flow_Lac"
$ws.Range("G39").Formula = '=AVERAGE(F$1:F$3)'
$ws.Range("H39").Formula = '=AVERAGE(F$4:F$6)'
$ws.Range("I39").Formula = '=AVERAGE(F$9:F$11)'
$ws.Range("J39").Formula = '=AVERAGE(F$13:F$16)'

# Row 40: standard error of the mean of column F per group.
$ws.Range("G40").Formula = '=STDEV(F$1:F$3)/SQRT(COUNT(F$1:F$3))'
$ws.Range("H40").Formula = '=STDEV(F$4:F$6)/SQRT(COUNT(F$4:F$6))'
$ws.Range("I40").Formula = '=STDEV(F$9:F$11)/SQRT(COUNT(F$9:F$11))'
$ws.Range("J40").Formula = '=STDEV(F$13:F$16)/SQRT(COUNT(F$13:F$16))'

# Scroll the view down to the newly added block and select it, matching
# the author's final on-screen selection.
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("F38:J40").Select()
